# Ajustes relatórios turmas A e C
# This script inserts a new column K (shifting the existing "Email" column
# from K to L) on the active sheet, sets the new K1 header to "2020-11-30"
# and fills the K2:K37 cells with the corresponding numeric scores taken
# from the diff. Rows whose K cell should remain empty are left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before K; this shifts the existing K column (which
# holds the "Email" header/values) one position to the right, into L,
# while preserving styles (including the bold/centered header style).
$ws.Columns("K").Insert()

# New header for the inserted column.
$ws.Range("K1").Value = "2020-11-30"

# New numeric values for K2:K37 (row number -> value). Rows not present
# in this table keep the empty cell left behind by the column insert.
$values = @{
    2  = 4.76
    3  = 80.95
    4  = 80
    5  = 48.57
    7  = 13.33
    8  = 11.43
    9  = 42.86
    10 = 28.57
    11 = 0
    13 = 82.86
    14 = 82.86
    17 = 11.43
    18 = 0
    19 = 82.86
    20 = 89.52
    21 = 87.62
    22 = 33.33
    23 = 0
    24 = 81.90000000000001
    27 = 0
    30 = 63.81
    31 = 85.70999999999999
    32 = 19.05
    33 = 86.67
    35 = 43.81
    36 = 80
    37 = 0
}

foreach ($row in $values.Keys) {
    $ws.Cells.Item($row, 11).Value = $values[$row]
}
